$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Journal de travail")
$ws2 = $wb.Worksheets.Item("Journal de bord")

# --- Sheet "Journal de travail": fill rows 4 and 5 ---
# Order of first-use of new strings must be: "Reçu sujet", "Divers",
# "Calculs moyennes, attente", "Prise en main d'un Fortinet..." so the
# shared-string table gets them appended in that exact order.

# Row 5 description first (introduces "Reçu sujet")
$ws1.Range("D5").Value = "Reçu sujet"
# Row 4 description (introduces "Divers")
$ws1.Range("D4").Value = "Divers"
# Row 4 remark (introduces "Calculs moyennes, attente")
$ws1.Range("E4").Value = "Calculs moyennes, attente"
# Row 5 remark (introduces Fortinet text)
$ws1.Range("E5").Value = "Prise en main d'un Fortinet: prise en main, comprendre le fonctionnement, implanter un/des firewall virtuels + début recherches Fortinet"

# Remaining cells for row 4
$ws1.Range("A4").Value = 44958
$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = "Documentation"
$ws1.Range("F4").Value = 0.63888888888888895
$ws1.Range("G4").Value = 0.66666666666666663

# Remaining cells for row 5
$ws1.Range("A5").Value = 44958
$ws1.Range("B5").Value = 1
$ws1.Range("C5").Value = "Documentation"
$ws1.Range("F5").Value = 0.66666666666666663
$ws1.Range("G5").Value = 0.70486111111111116

# Column width tweaks on "Journal de travail" (D & E grew wider)
$ws1.Columns.Item(4).ColumnWidth = 15.59
$ws1.Columns.Item(5).ColumnWidth = 81.25

# --- Sheet "Journal de bord": row 3 gets a new entry, rows 4-20 get the
#     date number format applied to column A (still blank) ---
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("A3").PasteSpecial(-4122) | Out-Null
$ws2.Range("A4:A20").PasteSpecial(-4122) | Out-Null
$ws2.Range("A3").Value = 44958
$ws2.Range("B3").Value = "Réception sujet Pré-TPI"
$ws2.Range("C3").Value = "Pas de cdc disponible"

# Column width tweak on "Journal de bord" (column B grew wider)
$ws2.Columns.Item(2).ColumnWidth = 13.92

# --- Selections: leave "Journal de bord" selection at A4, then finish on
#     "Journal de travail" (tab that stays active) at E27 ---
$ws2.Activate()
$ws2.Range("A4").Select() | Out-Null
$ws1.Activate()
$ws1.Range("E27").Select() | Out-Null
